$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 1495
$ws.Range("I94").Value = 1495
$ws.Range("K94").Value = 1495
$ws.Range("M94").Value = -1044
$ws.Range("H113").Value = 6949.25
$ws.Range("I113").Value = 5650
$ws.Range("J113").Value = 8248.5
$ws.Range("K113").Value = 5650
$ws.Range("L113").Value = 8248.5
$ws.Range("M113").Value = -2396
$ws.Range("N113").Value = -14756.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 23999.75
$ws.Range("I31").Value = 23999.75
$ws.Range("K31").Value = 23999.75
$ws.Range("M31").Value = -23705.75
$ws.Range("H61").Value = 2289.2144
$ws.Range("I61").Value = 2199.9092
$ws.Range("K61").Value = 2199.9092
$ws.Range("M61").Value = -1987.9092
$ws.Range("H74").Value = 844.13043
$ws.Range("I74").Value = 844.13043
$ws.Range("K74").Value = 844.13043
$ws.Range("M74").Value = 29.86956999999995
$ws.Range("H77").Value = 844.13043
$ws.Range("I77").Value = 844.13043
$ws.Range("K77").Value = 4220.65215
$ws.Range("M77").Value = 147.3478500000001
$ws.Range("H93").Value = 30298.666
$ws.Range("J93").Value = 35448
$ws.Range("L93").Value = 35448
$ws.Range("N93").Value = -40440
$ws.Range("H97").Value = 411.8889
$ws.Range("I97").Value = 402.125
$ws.Range("K97").Value = 402.125
$ws.Range("M97").Value = 93.875
$ws.Range("H132").Value = 2655.6538
$ws.Range("I132").Value = 2297.5454
$ws.Range("J132").Value = 4625.25
$ws.Range("K132").Value = 6892.6362
$ws.Range("L132").Value = 13875.75
$ws.Range("M132").Value = -4362.6362
$ws.Range("N132").Value = -18935.75
$ws.Range("H136").Value = 2289.2144
$ws.Range("I136").Value = 2199.9092
$ws.Range("K136").Value = 6599.7276
$ws.Range("M136").Value = -4049.7276

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7091
$ws.Range("I134").Value = 6762.8184
$ws.Range("J134").Value = 8294.333000000001
$ws.Range("K134").Value = 20288.4552
$ws.Range("L134").Value = 24882.999
$ws.Range("M134").Value = -17753.4552
$ws.Range("N134").Value = -29952.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1028.25
$ws.Range("I16").Value = 1269.4
$ws.Range("J16").Value = 626.3333
$ws.Range("K16").Value = 1269.4
$ws.Range("L16").Value = 626.3333
$ws.Range("M16").Value = -982.4000000000001
$ws.Range("N16").Value = -1200.3333
$ws.Range("H58").Value = 1919.1
$ws.Range("I58").Value = 1885.2858
$ws.Range("J58").Value = 1998
$ws.Range("K58").Value = 1885.2858
$ws.Range("L58").Value = 1998
$ws.Range("M58").Value = -1682.2858
$ws.Range("N58").Value = -2404
$ws.Range("H62").Value = 399.5
$ws.Range("I62").Value = 399.5
$ws.Range("K62").Value = 399.5
$ws.Range("M62").Value = 224.5
$ws.Range("H65").Value = 399.5
$ws.Range("I65").Value = 399.5
$ws.Range("K65").Value = 1997.5
$ws.Range("M65").Value = 1122.5
$ws.Range("H109").Value = 39642.5
$ws.Range("J109").Value = 39642.5
$ws.Range("L109").Value = 39642.5
$ws.Range("N109").Value = -41722.5
$ws.Range("H113").Value = 1028.25
$ws.Range("I113").Value = 1269.4
$ws.Range("J113").Value = 626.3333
$ws.Range("K113").Value = 1269.4
$ws.Range("L113").Value = 626.3333
$ws.Range("M113").Value = 900.5999999999999
$ws.Range("N113").Value = -4966.3333
$ws.Range("H122").Value = 1166.6666
$ws.Range("I122").Value = 1166.6666
$ws.Range("K122").Value = 3499.9998
$ws.Range("M122").Value = -1049.9998
$ws.Range("H132").Value = 4132.8
$ws.Range("I132").Value = 5600
$ws.Range("J132").Value = 3154.6667
$ws.Range("K132").Value = 16800
$ws.Range("L132").Value = 9464.000100000001
$ws.Range("M132").Value = -14270
$ws.Range("N132").Value = -14524.0001
$ws.Range("H134").Value = 1403.08
$ws.Range("I134").Value = 981.0454999999999
$ws.Range("J134").Value = 4498
$ws.Range("K134").Value = 2943.1365
$ws.Range("L134").Value = 13494
$ws.Range("M134").Value = -408.1364999999996
$ws.Range("N134").Value = -18564
$ws.Range("H136").Value = 1919.1
$ws.Range("I136").Value = 1885.2858
$ws.Range("J136").Value = 1998
$ws.Range("K136").Value = 5655.857400000001
$ws.Range("L136").Value = 5994
$ws.Range("M136").Value = -3105.857400000001
$ws.Range("N136").Value = -11094

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 963.2727
$ws.Range("I34").Value = 385.16666
$ws.Range("J34").Value = 1657
$ws.Range("K34").Value = 1155.49998
$ws.Range("L34").Value = 4971
$ws.Range("M34").Value = -1071.49998
$ws.Range("N34").Value = -5139
$ws.Range("H37").Value = 99998.664
$ws.Range("J37").Value = 99998.664
$ws.Range("L37").Value = 299995.992
$ws.Range("N37").Value = -300219.992
$ws.Range("H132").Value = 1229.875
$ws.Range("J132").Value = 1153.2727
$ws.Range("L132").Value = 10379.4543
$ws.Range("N132").Value = -15439.4543
$ws.Range("H137").Value = 569.8
$ws.Range("I137").Value = 569.8
$ws.Range("K137").Value = 1709.4
$ws.Range("M137").Value = 3390.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 8999
$ws.Range("J53").Value = 8999
$ws.Range("L53").Value = 8999
$ws.Range("N53").Value = -10261
$ws.Range("H122").Value = 3278.6
$ws.Range("I122").Value = 3050
$ws.Range("K122").Value = 9150
$ws.Range("M122").Value = -6700
$ws.Range("H126").Value = 6924.6
$ws.Range("I126").Value = 6331.6665
$ws.Range("K126").Value = 18994.9995
$ws.Range("M126").Value = -16524.9995
$ws.Range("H132").Value = 2202.2222
$ws.Range("I132").Value = 2004.1428
$ws.Range("K132").Value = 6012.428400000001
$ws.Range("M132").Value = -3482.428400000001
$ws.Range("H135").Value = 49999
$ws.Range("J135").Value = 49999
$ws.Range("L135").Value = 49999
$ws.Range("N135").Value = -60139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 36000
$ws.Range("J110").Value = 36000
$ws.Range("L110").Value = 36000
$ws.Range("N110").Value = -44180
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1057.0588
$ws.Range("I107").Value = 772.5
$ws.Range("K107").Value = 2317.5
$ws.Range("M107").Value = -397.5
$ws.Range("H122").Value = 3094.6
$ws.Range("J122").Value = 2898.6667
$ws.Range("L122").Value = 8696.000100000001
$ws.Range("N122").Value = -13596.0001
$ws.Range("H123").Value = 297500
$ws.Range("J123").Value = 297500
$ws.Range("L123").Value = 297500
$ws.Range("N123").Value = -307300
$ws.Range("H125").Value = 118000
$ws.Range("J125").Value = 118000
$ws.Range("L125").Value = 118000
$ws.Range("N125").Value = -127840
$ws.Range("H132").Value = 2669.6667
$ws.Range("I132").Value = 2004
$ws.Range("K132").Value = 6012
$ws.Range("M132").Value = -3482
